# Add a "Save" column (H) to the s_vals sheet, matching the header
# formatting already used by the other header cells (B1:G1), and a
# numeric value of 0 for the data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell in H1, same style as the existing headers.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# New data cell in H2 holding the numeric value 0.
$ws.Range("H2").Value = 0
